# GDP by Industries USA - IoP by Industries UK
# Fix header typo/casing in the "Assets _M" column header and
# leave the selection where the author last clicked (H12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1 holds the header "Assets _M" -> correct to "assets _m"
$ws.Range("C1").Value = "assets _m"

# Move/restore the active selection to H12 (was A34:D34)
$ws.Range("H12").Select()
